# MBD Problem Solver Input.xlsx - update input values on the "Input" sheet
# (the active/tab-selected worksheet) plus a couple of small formatting/
# selection tweaks, matching the author's latest edits to the tool's
# default example values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Torsion Variables header (Q1): re-apply the same header formatting
# used by the neighboring group headers (L1/M1 "C-S Dimensions" /
# "Bending Variables") so the column reads as its own boxed header cell.
$ws.Range("L1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value = "Torsion Variables"

# --- Row 2 (example/default values) ---
# Torsion shear modulus variable (L2) and Cross-section target radius (Q2)
$ws.Range("L2").Value = 0.024
$ws.Range("Q2").Value = 0.024

# Cross-section total beam length (O2)
$ws.Range("O2").Value = 300

# --- Row 3 (example/default values) ---
# Bending: target distance from neutral axis (M3)
$ws.Range("M3").Value = 0

# Cross-section dimension (Q3)
$ws.Range("Q3").Value = 0.9

# --- Cursor/selection, matching where the user left off editing ---
$ws.Range("N3").Select() | Out-Null
